# Natmi following Dr Hou advice
#
# Updates the LR-pair (Fgf1-Fgfr3) results table on the active sheet:
#   - Rows 2-7 get recomputed metric values (ligand-expressing-cell counts
#     changed from 1 to 3, which ripples through the derived columns), and
#     the "Target cluster" (column D) / "Sending cluster" (column A) labels
#     are reassigned for some rows.
#   - Three new rows (8, 9, 10) are appended for the "sCs" sending cluster
#     against each of the three target clusters (ECs, FAPs, sCs), growing
#     the used range from A1:T7 to A1:T10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One entry per data row: RowNum is the worksheet row, Col<L> holds the
# value for column <L> (A..T), matching the header layout:
#   A Sending cluster            K Receptor-expressing cells
#   B Ligand symbol              L Receptor detection rate
#   C Receptor symbol            M Receptor average expression value
#   D Target cluster             N Receptor total expression value
#   E Ligand-expressing cells    O Receptor derived specificity (avg)
#   F Ligand detection rate      P Receptor derived specificity (total)
#   G Ligand average expr value  Q Edge average expression weight
#   H Ligand total expr value    R Edge total expression weight
#   I Ligand derived spec (avg)  S Edge average expression derived spec
#   J Ligand derived spec (tot)  T Edge total expression derived spec
$rowData = @(
    @{ RowNum=2;  ColA="ECs";  ColB="Fgf1"; ColC="Fgfr3"; ColD="ECs";
       ColE=3; ColF=1; ColG=0.6873773333333334; ColH=2.062132;
       ColI=0.02660947569874856; ColJ=0.02660947569874856;
       ColK=3; ColL=1; ColM=1.656156333333334; ColN=4.968469000000001;
       ColO=0.6151212440816572; ColP=0.6151212440816572;
       ColQ=1.138404323989778; ColR=10.245638915908;
       ColS=0.01636805379617484; ColT=0.01636805379617484 },

    @{ RowNum=3;  ColA="ECs";  ColB="Fgf1"; ColC="Fgfr3"; ColD="FAPs";
       ColE=3; ColF=1; ColG=0.6873773333333334; ColH=2.062132;
       ColI=0.02660947569874856; ColJ=0.02660947569874856;
       ColK=2; ColL=0.6666666666666666; ColM=0.5648773333333333; ColN=1.694632;
       ColO=0.2098038941373262; ColP=0.2098038941373262;
       ColQ=0.3882838750471111; ColR=3.494554875424;
       ColS=0.005582771622549997; ColT=0.005582771622549996 },

    @{ RowNum=4;  ColA="ECs";  ColB="Fgf1"; ColC="Fgfr3"; ColD="sCs";
       ColE=3; ColF=1; ColG=0.6873773333333334; ColH=2.062132;
       ColI=0.02660947569874856; ColJ=0.02660947569874856;
       ColK=3; ColL=1; ColM=0.4713726666666667; ColN=1.414118;
       ColO=0.1750748617810164; ColP=0.1750748617810165;
       ColQ=0.3240108866195556; ColR=2.916097979576;
       ColS=0.00465865028002372; ColT=0.004658650280023721 },

    @{ RowNum=5;  ColA="FAPs"; ColB="Fgf1"; ColC="Fgfr3"; ColD="ECs";
       ColE=3; ColF=1; ColG=5.913984666666667; ColH=17.741954;
       ColI=0.2289398029860915; ColJ=0.2289398029860915;
       ColK=3; ColL=1; ColM=1.656156333333334; ColN=4.968469000000001;
       ColO=0.6151212440816572; ColP=0.6151212440816572;
       ColQ=9.794483160936224; ColR=88.15034844842602;
       ColS=0.1408257364326141; ColT=0.1408257364326141 },

    @{ RowNum=6;  ColA="FAPs"; ColB="Fgf1"; ColC="Fgfr3"; ColD="FAPs";
       ColE=3; ColF=1; ColG=5.913984666666667; ColH=17.741954;
       ColI=0.2289398029860915; ColJ=0.2289398029860915;
       ColK=2; ColL=0.6666666666666666; ColM=0.5648773333333333; ColN=1.694632;
       ColO=0.2098038941373262; ColP=0.2098038941373262;
       ColQ=3.340675887880889; ColR=30.066082990928;
       ColS=0.04803246218951425; ColT=0.04803246218951425 },

    @{ RowNum=7;  ColA="FAPs"; ColB="Fgf1"; ColC="Fgfr3"; ColD="sCs";
       ColE=3; ColF=1; ColG=5.913984666666667; ColH=17.741954;
       ColI=0.2289398029860915; ColJ=0.2289398029860915;
       ColK=3; ColL=1; ColM=0.4713726666666667; ColN=1.414118;
       ColO=0.1750748617810164; ColP=0.1750748617810165;
       ColQ=2.787690722952445; ColR=25.089216506572;
       ColS=0.0400816043639631; ColT=0.04008160436396311 },

    @{ RowNum=8;  ColA="sCs";  ColB="Fgf1"; ColC="Fgfr3"; ColD="ECs";
       ColE=3; ColF=1; ColG=19.230689; ColH=57.692067;
       ColI=0.7444507213151601; ColJ=0.7444507213151601;
       ColK=3; ColL=1; ColM=1.656156333333334; ColN=4.968469000000001;
       ColO=0.6151212440816572; ColP=0.6151212440816572;
       ColQ=31.84902738171368; ColR=286.6412464354231;
       ColS=0.4579274538528684; ColT=0.4579274538528684 },

    @{ RowNum=9;  ColA="sCs";  ColB="Fgf1"; ColC="Fgfr3"; ColD="FAPs";
       ColE=3; ColF=1; ColG=19.230689; ColH=57.692067;
       ColI=0.7444507213151601; ColJ=0.7444507213151601;
       ColK=2; ColL=0.6666666666666666; ColM=0.5648773333333333; ColN=1.694632;
       ColO=0.2098038941373262; ColP=0.2098038941373262;
       ColQ=10.86298032048267; ColR=97.76682288434399;
       ColS=0.156188660325262; ColT=0.156188660325262 },

    @{ RowNum=10; ColA="sCs";  ColB="Fgf1"; ColC="Fgfr3"; ColD="sCs";
       ColE=3; ColF=1; ColG=19.230689; ColH=57.692067;
       ColI=0.7444507213151601; ColJ=0.7444507213151601;
       ColK=3; ColL=1; ColM=0.4713726666666667; ColN=1.414118;
       ColO=0.1750748617810164; ColP=0.1750748617810165;
       ColQ=9.064821155767335; ColR=81.58339040190602;
       ColS=0.1303346071370297; ColT=0.1303346071370297 }
)

# NOTE: the hashtable key for the worksheet row number is "RowNum", not "R" —
# PowerShell hashtable keys are case-insensitive, so a key named "R" would
# collide with column R's "ColR" only by coincidence of naming; RowNum avoids
# any ambiguity with the column-letter keys entirely.
$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rowData) {
    $targetRow = $row.RowNum
    for ($i = 0; $i -lt $colOrder.Count; $i++) {
        $colLetter = $colOrder[$i]
        $colIndex = $i + 1
        $key = "Col" + $colLetter
        $ws.Cells.Item($targetRow, $colIndex).Value = $row[$key]
    }
}

Write-Output "Updated rows 2-10 (added rows 8-10) on $($ws.Name)"
